# Update "想去人数" (interested-count) figures to the latest scraped values.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 5571
$ws1.Range("F4").Value = 640
$ws1.Range("F6").Value = 832
$ws1.Range("F7").Value = 50
$ws1.Range("F8").Value = 368
$ws1.Range("F9").Value = 3

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 19

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 5571
$ws4.Range("F4").Value = 640
$ws4.Range("F6").Value = 832
$ws4.Range("F7").Value = 50
$ws4.Range("F9").Value = 368
$ws4.Range("F10").Value = 3
$ws4.Range("F13").Value = 19

$wb.Save()
